$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert a new row at position 9 (shifts old rows 9-21 down to 10-22)
$ws.Rows.Item(9).Insert()

# Copy the formatting of the whole row 7 into row 9
$ws.Rows.Item(7).Copy()
$ws.Rows.Item(9).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "done"
